$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3..12, column F ("想去人数")
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 532
$wsExhibit.Range("F4").Value = 174
$wsExhibit.Range("F5").Value = 249
$wsExhibit.Range("F6").Value = 380
$wsExhibit.Range("F7").Value = 236
$wsExhibit.Range("F8").Value = 2275
$wsExhibit.Range("F9").Value = 384
$wsExhibit.Range("F10").Value = 5639
$wsExhibit.Range("F11").Value = 135
$wsExhibit.Range("F12").Value = 369

# Sheet "全部类型" (All Types) - rows 4..15, column F ("想去人数")
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 532
$wsAll.Range("F5").Value = 174
$wsAll.Range("F6").Value = 249
$wsAll.Range("F7").Value = 380
$wsAll.Range("F8").Value = 236
$wsAll.Range("F11").Value = 2275
$wsAll.Range("F12").Value = 384
$wsAll.Range("F13").Value = 5639
$wsAll.Range("F14").Value = 135
$wsAll.Range("F15").Value = 369
